# Updates cryptos list values (price/volume columns, plus two swapped
# coin rows) to match the refreshed GitHub Actions data pull.
#
# Price-like text in column D (e.g. "56.045.76") uses dots as thousands
# separators and must stay a literal text string, not be reinterpreted as
# a number. Assigning via Value2 with a leading literal apostrophe forces
# text entry (mirrors typing '3.34 in Excel); resetting Style to "Normal"
# afterwards drops the transient quote-prefix number format so no stray
# style gets attached to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.Value2 = "'" + $text
    $cell.Style = 'Normal'
}

# Row 2
Set-TextCell 'D2' '56.271.79'
Set-TextCell 'E2' '  +9.67%  '
# Row 3
Set-TextCell 'D3' '3.224.95'
Set-TextCell 'E3' '  +4.19%  '
# Row 4
Set-TextCell 'E4' '  -0.01%  '
# Row 5
Set-TextCell 'D5' '397.29'
Set-TextCell 'E5' '  +2.40%  '
# Row 6
Set-TextCell 'D6' '111.02'
Set-TextCell 'E6' '  +7.31%  '
# Row 7
Set-TextCell 'D7' '0.554'
Set-TextCell 'E7' '  +3.05%  '
# Row 8
Set-TextCell 'D8' '0.999'
Set-TextCell 'E8' '  -0.04%  '
# Row 9
Set-TextCell 'D9' '0.618'
Set-TextCell 'E9' '  +5.76%  '
# Row 10
Set-TextCell 'D10' '39.23'
Set-TextCell 'E10' '  +6.18%  '
# Row 11
Set-TextCell 'E11' '  +7.02%  '
# Row 12
Set-TextCell 'E12' '  +2.03%  '
# Row 13
Set-TextCell 'D13' '3.732.29'
Set-TextCell 'E13' '  +4.34%  '
# Row 14
Set-TextCell 'D14' '8.08'
Set-TextCell 'E14' '  +4.43%  '
# Row 15
Set-TextCell 'D15' '19.04'
Set-TextCell 'E15' '  +3.03%  '
# Row 16
Set-TextCell 'D16' '3.230.10'
Set-TextCell 'E16' '  +4.35%  '
# Row 17
Set-TextCell 'E17' '  +5.15%  '
# Row 18
Set-TextCell 'D18' '10.85'
Set-TextCell 'E18' '  +2.00%  '
# Row 19
Set-TextCell 'D19' '56.046.77'
Set-TextCell 'E19' '  +9.04%  '
# Row 20
Set-TextCell 'D20' '3.34'
Set-TextCell 'E20' '  +3.31%  '
# Row 21
Set-TextCell 'E21' '  +6.55%  '
# Row 22
Set-TextCell 'D22' '12.97'
Set-TextCell 'E22' '  +3.75%  '
# Row 23
Set-TextCell 'D23' '298.40'
Set-TextCell 'E23' '  +12.30%  '
# Row 24
Set-TextCell 'D24' '75.55'
Set-TextCell 'E24' '  +8.00%  '
# Row 25
Set-TextCell 'E25' '  +1.74%  '
# Row 26
Set-TextCell 'D26' '8.18'
Set-TextCell 'E26' '  +2.48%  '
# Row 27
Set-TextCell 'D27' '28.15'
Set-TextCell 'E27' '  +3.06%  '
# Row 28
Set-TextCell 'D28' '7.47'
Set-TextCell 'E28' '  +3.47%  '
# Row 29
Set-TextCell 'D29' '0.173'
Set-TextCell 'E29' '  +4.61%  '
# Row 30
Set-TextCell 'E30' '  +0.39%  '
# Row 31
Set-TextCell 'E31' '  +4.23%  '
# Row 32
Set-TextCell 'D32' '11.14'
Set-TextCell 'E32' '  +7.23%  '
# Row 33
Set-TextCell 'D33' '0.0497'
Set-TextCell 'E33' '  +4.75%  '
# Row 34
Set-TextCell 'D34' '36.43'
Set-TextCell 'E34' '  +1.26%  '
# Row 35
Set-TextCell 'D35' '2.13'
Set-TextCell 'E35' '  +2.71%  '
# Row 36
Set-TextCell 'D36' '51.38'
Set-TextCell 'E36' '  +3.17%  '
# Row 37
Set-TextCell 'D37' '3.14'
Set-TextCell 'E37' '  +25.69%  '
# Row 38
Set-TextCell 'D38' '3.54'
Set-TextCell 'E38' '  +4.49%  '
# Row 39
Set-TextCell 'D39' '0.999'
Set-TextCell 'E39' '  -0.05%  '
# Row 40
Set-TextCell 'D40' '137.58'
Set-TextCell 'E40' '  +5.35%  '
# Row 41
Set-TextCell 'D41' '17.38'
Set-TextCell 'E41' '  +5.04%  '
# Row 42
Set-TextCell 'D42' '1.92'
Set-TextCell 'E42' '  +3.30%  '
# Row 43
Set-TextCell 'D43' '4.00'
Set-TextCell 'E43' '  +3.56%  '
# Row 44
Set-TextCell 'E44' '  +3.35%  '
# Row 45
Set-TextCell 'D45' '0.285'
Set-TextCell 'E45' '  -1.97%  '
# Row 46
Set-TextCell 'D46' '22.15'
Set-TextCell 'E46' '  +0.67%  '
# Row 47
Set-TextCell 'D47' '2.20'
Set-TextCell 'E47' '  +54.44%  '
# Row 48
Set-TextCell 'B48' 'ApeXProtocol'
Set-TextCell 'C48' 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextCell 'D48' '2.47'
Set-TextCell 'E48' '  -1.97%  '
# Row 49
Set-TextCell 'B49' 'WEMIXToken'
Set-TextCell 'C49' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextCell 'D49' '2.09'
Set-TextCell 'E49' '  +0.02%  '
# Row 50
Set-TextCell 'D50' '2.120.13'
Set-TextCell 'E50' '  +2.24%  '
# Row 51
Set-TextCell 'E51' '  +9.50%  '

Write-Output "Updated $([int]94) cells across the cryptos list"
